$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6 (H6 anchor old=777.9167)
$ws.Range("H6").Value = 1105.5
$ws.Range("I6").Value = 1198.8
$ws.Range("J6").Value = 950
$ws.Range("K6").Value = 3596.4
$ws.Range("L6").Value = 2850
$ws.Range("M6").Value = -3484.4
$ws.Range("N6").Value = -3074

# Row 33 (H33 anchor old=262.78262)
$ws.Range("H33").Value = 298.1111
$ws.Range("I33").Value = 249.4
$ws.Range("K33").Value = 249.4
$ws.Range("M33").Value = -20.40000000000001

# Row 62 (H62 anchor old=4000)
$ws.Range("H62").Value = 3933.3333
$ws.Range("I62").Value = 2900
$ws.Range("K62").Value = 2900
$ws.Range("M62").Value = -2276

# Row 65 (H65 anchor old=4000)
$ws.Range("H65").Value = 3933.3333
$ws.Range("I65").Value = 2900
$ws.Range("K65").Value = 14500
$ws.Range("M65").Value = -11380

# Row 107 (H107 anchor old=185)
$ws.Range("H107").Value = 172.66667
$ws.Range("I107").Value = 213.5
$ws.Range("K107").Value = 213.5
$ws.Range("M107").Value = 1706.5

# Row 132 (H132 anchor old=1272.0193)
$ws.Range("H132").Value = 1250.9434
$ws.Range("I132").Value = 1290.3334
$ws.Range("K132").Value = 3871.0002
$ws.Range("M132").Value = -1341.0002

# Row 141 (H141 anchor old=6999)
$ws.Range("H141").Value = 6134.857
$ws.Range("I141").Value = 3236.25
$ws.Range("K141").Value = 9708.75
$ws.Range("M141").Value = -4528.75


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (H2 anchor old=774.6)
$ws.Range("H2").Value = 739
$ws.Range("I2").Value = 548.75
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 548.75
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -435.75
$ws.Range("N2").Value = -1726

# Row 11 (H11 anchor old=3133)
$ws.Range("H11").Value = 11500000
$ws.Range("I11").Value = 11500000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 11500000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -11499856
$ws.Range("N11").ClearContents()

# Row 40 (H40 anchor old=50000)
$ws.Range("H40").Value = 49500
$ws.Range("J40").Value = 49500
$ws.Range("L40").Value = 49500
$ws.Range("N40").Value = -49852

# Row 61 (H61 anchor old=2081.2307)
$ws.Range("H61").Value = 1986.0714
$ws.Range("I61").Value = 1632.5
$ws.Range("K61").Value = 1632.5
$ws.Range("M61").Value = -1420.5

# Row 74 (H74 anchor old=2380.6428)
$ws.Range("H74").Value = 2288.2666
$ws.Range("J74").Value = 2818.4285
$ws.Range("L74").Value = 2818.4285
$ws.Range("N74").Value = -4566.4285

# Row 77 (H77 anchor old=2380.6428)
$ws.Range("H77").Value = 2288.2666
$ws.Range("J77").Value = 2818.4285
$ws.Range("L77").Value = 14092.1425
$ws.Range("N77").Value = -22828.1425

# Row 116 (H116 anchor old=774.6)
$ws.Range("H116").Value = 739
$ws.Range("I116").Value = 548.75
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 548.75
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1745.25
$ws.Range("N116").Value = -6088

# Row 130 (H130 anchor old=43297.332)
$ws.Range("H130").Value = 43296
$ws.Range("J130").Value = 43296
$ws.Range("L130").Value = 43296
$ws.Range("N130").Value = -53336

# Row 132 (H132 anchor old=1459.1395)
$ws.Range("H132").Value = 1385.3112
$ws.Range("I132").Value = 1118.3684
$ws.Range("K132").Value = 3355.1052
$ws.Range("M132").Value = -825.1052

# Row 136 (H136 anchor old=2081.2307)
$ws.Range("H136").Value = 1986.0714
$ws.Range("I136").Value = 1632.5
$ws.Range("K136").Value = 4897.5
$ws.Range("M136").Value = -2347.5


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (H3 anchor old=774.6)
$ws.Range("H3").Value = 739
$ws.Range("I3").Value = 548.75
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 548.75
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -434.75
$ws.Range("N3").Value = -1728

# Row 99 (H99 anchor old=64561)
$ws.Range("H99").Value = 57498.055
$ws.Range("I99").Value = 126129.5
$ws.Range("K99").Value = 126129.5
$ws.Range("M99").Value = -124631.5


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4 (H4 anchor old=8675309)
$ws.Range("H4").Value = 6508982
$ws.Range("J4").Value = 6508982
$ws.Range("L4").Value = 6508982
$ws.Range("N4").Value = -6509206

# Row 25 (H25 anchor old=4255)
$ws.Range("H25").Value = 5604
$ws.Range("I25").Value = 5604
$ws.Range("K25").Value = 5604
$ws.Range("M25").Value = -5430

# Row 94 (H94 anchor old=7966.2)
$ws.Range("H94").Value = 625.53845
$ws.Range("I94").Value = 660.125
$ws.Range("K94").Value = 660.125
$ws.Range("M94").Value = -209.125

# Row 127 (H127 anchor old=0)
$ws.Range("H127").Value = 33000
$ws.Range("I127").Value = 33000
$ws.Range("K127").Value = 33000
$ws.Range("M127").Value = -28040

# Row 134 (H134 anchor old=4470.909)
$ws.Range("H134").Value = 4306.6665
$ws.Range("I134").Value = 4518
$ws.Range("K134").Value = 13554
$ws.Range("M134").Value = -11019


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 9 (H9 anchor old=2500352.8)
$ws.Range("H9").Value = 3336833
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 30000
$ws.Range("N9").Value = -30448

# Row 13 (H13 anchor old=334.33334)
$ws.Range("H13").Value = 666.6667
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -1332
$ws.Range("N13").Value = -3336

# Row 82 (H82 anchor old=7336)
$ws.Range("H82").Value = 6856
$ws.Range("I82").Value = 7776.5
$ws.Range("J82").Value = 5015
$ws.Range("K82").Value = 23329.5
$ws.Range("L82").Value = 15045
$ws.Range("M82").Value = -22923.5
$ws.Range("N82").Value = -15857

# Row 85 (H85 anchor old=7336)
$ws.Range("H85").Value = 6856
$ws.Range("I85").Value = 7776.5
$ws.Range("J85").Value = 5015
$ws.Range("K85").Value = 23329.5
$ws.Range("L85").Value = 15045
$ws.Range("M85").Value = -21925.5
$ws.Range("N85").Value = -17853

# Row 92 (H92 anchor old=220.5)
$ws.Range("H92").Value = 220.33333
$ws.Range("I92").Value = 231.8
$ws.Range("J92").Value = 163
$ws.Range("K92").Value = 695.4000000000001
$ws.Range("L92").Value = 489
$ws.Range("M92").Value = 552.5999999999999
$ws.Range("N92").Value = -2985

# Row 132 (H132 anchor old=5993.2354)
$ws.Range("H132").Value = 5529.615
$ws.Range("J132").Value = 8786.857
$ws.Range("L132").Value = 79081.713
$ws.Range("N132").Value = -84141.713


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5 (H5 anchor old=50)
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 113 (H113 anchor old=1757.5)
$ws.Range("H113").Value = 1509.9
$ws.Range("I113").Value = 1137.375
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1137.375
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1032.625
$ws.Range("N113").Value = -7340

# Row 123 (H123 anchor old=19485)
$ws.Range("H123").Value = 19483.334
$ws.Range("J123").Value = 19483.334
$ws.Range("L123").Value = 19483.334
$ws.Range("N123").Value = -24383.334


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (H7 anchor old=5002.8)
$ws.Range("H7").Value = 4521.385
$ws.Range("I7").Value = 3350.5
$ws.Range("K7").Value = 3350.5
$ws.Range("M7").Value = -3238.5

# Row 46 (H46 anchor old=3183.7188)
$ws.Range("H46").Value = 3490.3704
$ws.Range("I46").Value = 2750
$ws.Range("J46").Value = 3582.9167
$ws.Range("K46").Value = 2750
$ws.Range("L46").Value = 3582.9167
$ws.Range("M46").Value = -2562
$ws.Range("N46").Value = -3958.9167

# Row 61 (H61 anchor old=2623.6667)
$ws.Range("H61").Value = 2186
$ws.Range("I61").Value = 915
$ws.Range("K61").Value = 915
$ws.Range("M61").Value = -713

# Row 82 (H82 anchor old=2958.8)
$ws.Range("H82").Value = 3000
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2639
$ws.Range("N82").ClearContents()

# Row 85 (H85 anchor old=2958.8)
$ws.Range("H85").Value = 3000
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1752
$ws.Range("N85").ClearContents()

# Row 113 (H113 anchor old=2623.6667)
$ws.Range("H113").Value = 2186
$ws.Range("I113").Value = 915
$ws.Range("K113").Value = 915
$ws.Range("M113").Value = 1255

# Row 126 (H126 anchor old=5002.8)
$ws.Range("H126").Value = 4521.385
$ws.Range("I126").Value = 3350.5
$ws.Range("K126").Value = 10051.5
$ws.Range("M126").Value = -7581.5


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122 (H122 anchor old=1711.8334)
$ws.Range("H122").Value = 1666.24
$ws.Range("J122").Value = 2712
$ws.Range("L122").Value = 8136
$ws.Range("N122").Value = -13036

# Row 132 (H132 anchor old=1137.0652)
$ws.Range("H132").Value = 1145.6666
$ws.Range("I132").Value = 902.5263
$ws.Range("K132").Value = 2707.5789
$ws.Range("M132").Value = -177.5789

# Row 136 (H136 anchor old=1781.8572)
$ws.Range("H136").Value = 1695.7333
$ws.Range("I136").Value = 1516
$ws.Range("K136").Value = 4548
$ws.Range("M136").Value = -1998

